# Scheduled-runner update: refresh currentAveragePrice/Leve-profit columns
# (H, I, J, K, L, M, N) for a set of rows across the ALC/ARM/BSM/CRP/CUL/
# GSM/LTW/WVR sheets with newly pulled market-board data. A few rows also
# gain/lose a trailing profit cell (M/N) that previously did not exist.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 208.4375
$ws.Range("I9").Value = 80.454544
$ws.Range("K9").Value = 80.454544
$ws.Range("M9").Value = 88.545456

$ws.Range("H11").Value = 194502.5
$ws.Range("I11").Value = 194502.5
$ws.Range("K11").Value = 194502.5
$ws.Range("M11").Value = -194362.5

$ws.Range("H53").Value = 439.3846
$ws.Range("I53").Value = 82.666664
$ws.Range("J53").Value = 745.1429000000001
$ws.Range("K53").Value = 82.666664
$ws.Range("L53").Value = 745.1429000000001
$ws.Range("M53").Value = 554.333336
$ws.Range("N53").Value = -2019.1429

$ws.Range("H62").Value = 17831.834
$ws.Range("I62").Value = 12398.2
$ws.Range("K62").Value = 12398.2
$ws.Range("M62").Value = -11774.2

$ws.Range("H65").Value = 17831.834
$ws.Range("I65").Value = 12398.2
$ws.Range("K65").Value = 61991
$ws.Range("M65").Value = -58871

$ws.Range("H132").Value = 16766
$ws.Range("I132").Value = 18576.154
$ws.Range("K132").Value = 55728.462
$ws.Range("M132").Value = -53198.462

$ws.Range("H137").Value = 26319166
$ws.Range("I137").Value = 55556610
$ws.Range("K137").Value = 166669830
$ws.Range("M137").Value = -166667280

$ws.Range("H138").Value = 3214.5
$ws.Range("J138").Value = 3526.1936
$ws.Range("L138").Value = 10578.5808
$ws.Range("N138").Value = -20858.5808

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 175538.72
$ws.Range("I32").Value = 305151.12
$ws.Range("K32").Value = 305151.12
$ws.Range("M32").Value = -304864.12

$ws.Range("H45").Value = 1900
$ws.Range("I45").Value = 1750
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 1750
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -1373
$ws.Range("N45").Value = -2754

$ws.Range("H74").Value = 1295775.9
$ws.Range("I74").Value = 1986412
$ws.Range("J74").Value = 6588.2666
$ws.Range("K74").Value = 1986412
$ws.Range("L74").Value = 6588.2666
$ws.Range("M74").Value = -1985538
$ws.Range("N74").Value = -8336.266599999999

$ws.Range("H77").Value = 1295775.9
$ws.Range("I77").Value = 1986412
$ws.Range("J77").Value = 6588.2666
$ws.Range("K77").Value = 9932060
$ws.Range("L77").Value = 32941.333
$ws.Range("M77").Value = -9927692
$ws.Range("N77").Value = -41677.333

$ws.Range("H122").Value = 1414.1
$ws.Range("I122").Value = 987.2857
$ws.Range("K122").Value = 2961.8571
$ws.Range("M122").Value = -511.8571000000002

$ws.Range("H132").Value = 696796.9399999999
$ws.Range("I132").Value = 808702.4399999999
$ws.Range("J132").Value = 2983
$ws.Range("K132").Value = 2426107.32
$ws.Range("L132").Value = 8949
$ws.Range("M132").Value = -2423577.32
$ws.Range("N132").Value = -14009

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1655.4
$ws.Range("I5").Value = 588.5
$ws.Range("J5").Value = 2366.6667
$ws.Range("K5").Value = 588.5
$ws.Range("L5").Value = 2366.6667
$ws.Range("M5").Value = -475.5
$ws.Range("N5").Value = -2592.6667

$ws.Range("H20").Value = 45711.543
$ws.Range("I20").Value = 57214.58
$ws.Range("J20").Value = 2000
$ws.Range("K20").Value = 57214.58
$ws.Range("L20").Value = 2000
$ws.Range("M20").Value = -56967.58
$ws.Range("N20").Value = -2494

$ws.Range("H133").Value = 101665.664
$ws.Range("J133").Value = 101665.664
$ws.Range("L133").Value = 101665.664
$ws.Range("N133").Value = -111785.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 90963.09
$ws.Range("I7").Value = 111157.78
$ws.Range("K7").Value = 111157.78
$ws.Range("M7").Value = -111044.78

$ws.Range("H8").Value = 259.5
$ws.Range("I8").Value = 9
$ws.Range("J8").Value = 510
$ws.Range("K8").Value = 9
$ws.Range("L8").Value = 510
$ws.Range("M8").Value = 131
$ws.Range("N8").Value = -790

$ws.Range("H15").Value = 1842.6666
$ws.Range("I15").Value = 2754
$ws.Range("J15").Value = 20
$ws.Range("K15").Value = 2754
$ws.Range("L15").Value = 20
$ws.Range("M15").Value = -2584
$ws.Range("N15").Value = -360

$ws.Range("H31").Value = 1986775.9
$ws.Range("I31").Value = 2418195.5
$ws.Range("K31").Value = 2418195.5
$ws.Range("M31").Value = -2417900.5

$ws.Range("H34").Value = 1986775.9
$ws.Range("I34").Value = 2418195.5
$ws.Range("K34").Value = 2418195.5
$ws.Range("M34").Value = -2417993.5

$ws.Range("H44").Value = 12500.5
$ws.Range("I44").Value = 10000
$ws.Range("J44").Value = 15001
$ws.Range("K44").Value = 10000
$ws.Range("L44").Value = 15001
$ws.Range("M44").Value = -9558
$ws.Range("N44").Value = -15885

$ws.Range("H93").Value = 9132.333000000001
$ws.Range("I93").Value = 9132.333000000001
$ws.Range("K93").Value = 9132.333000000001
$ws.Range("M93").Value = -7260.333000000001

$ws.Range("H99").Value = 25944.309
$ws.Range("I99").Value = 25549.9
$ws.Range("K99").Value = 25549.9
$ws.Range("M99").Value = -24051.9

$ws.Range("H110").Value = 30638
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws.Range("H126").Value = 25944.309
$ws.Range("I126").Value = 25549.9
$ws.Range("K126").Value = 76649.70000000001
$ws.Range("M126").Value = -74179.70000000001

$ws.Range("H141").Value = 171517
$ws.Range("J141").Value = 203563.83
$ws.Range("L141").Value = 203563.83
$ws.Range("N141").Value = -213923.83

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 8034.8276
$ws.Range("J68").Value = 8034.8276
$ws.Range("L68").Value = 24104.4828
$ws.Range("N68").Value = -25726.4828

$ws.Range("H69").Value = 6714.143
$ws.Range("I69").Value = 999.5
$ws.Range("J69").Value = 9000
$ws.Range("K69").Value = 2998.5
$ws.Range("L69").Value = 27000
$ws.Range("M69").Value = -2187.5
$ws.Range("N69").Value = -28622

$ws.Range("H71").Value = 8034.8276
$ws.Range("J71").Value = 8034.8276
$ws.Range("L71").Value = 72313.44839999999
$ws.Range("N71").Value = -80425.44839999999

$ws.Range("H72").Value = 6714.143
$ws.Range("I72").Value = 999.5
$ws.Range("J72").Value = 9000
$ws.Range("K72").Value = 8995.5
$ws.Range("L72").Value = 81000
$ws.Range("M72").Value = -4939.5
$ws.Range("N72").Value = -89112

$ws.Range("H97").Value = 317.94623
$ws.Range("J97").Value = 566.25
$ws.Range("L97").Value = 1698.75
$ws.Range("N97").Value = -2690.75

$ws.Range("H107").Value = 3731.2
$ws.Range("J107").Value = 5307
$ws.Range("L107").Value = 15921
$ws.Range("N107").Value = -19761

$ws.Range("H109").Value = 15000
$ws.Range("J109").Value = 15000
$ws.Range("L109").Value = 45000
$ws.Range("N109").Value = -47080

$ws.Range("H113").Value = 1405.7693
$ws.Range("I113").Value = 692.3333
$ws.Range("J113").Value = 1619.8
$ws.Range("K113").Value = 2076.9999
$ws.Range("L113").Value = 4859.4
$ws.Range("M113").Value = 93.0001000000002
$ws.Range("N113").Value = -9199.4

$ws.Range("H131").Value = 5604.364
$ws.Range("J131").Value = 8164.2856
$ws.Range("L131").Value = 24492.8568
$ws.Range("N131").Value = -34572.8568

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 14710.4
$ws.Range("I80").Value = 9332.333000000001
$ws.Range("K80").Value = 9332.333000000001
$ws.Range("M80").Value = -8334.333000000001

$ws.Range("H83").Value = 14710.4
$ws.Range("I83").Value = 9332.333000000001
$ws.Range("K83").Value = 46661.665
$ws.Range("M83").Value = -41669.665

$ws.Range("H102").Value = 2088.2222
$ws.Range("I102").Value = 2034.5883
$ws.Range("K102").Value = 2034.5883
$ws.Range("M102").Value = -412.5882999999999

$ws.Range("H107").Value = 4910.037
$ws.Range("J107").Value = 2367.1428
$ws.Range("L107").Value = 2367.1428
$ws.Range("N107").Value = -6207.1428

$ws.Range("H122").Value = 86014.92
$ws.Range("I122").Value = 133899.62
$ws.Range("K122").Value = 401698.86
$ws.Range("M122").Value = -399248.86

$ws.Range("H132").Value = 5155.9077
$ws.Range("J132").Value = 7457.643
$ws.Range("L132").Value = 22372.929
$ws.Range("N132").Value = -27432.929

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3231.3333
$ws.Range("J46").Value = 7102.6
$ws.Range("L46").Value = 7102.6
$ws.Range("N46").Value = -7478.6

$ws.Range("H68").Value = 1833.3334
$ws.Range("I68").Value = 1833.3334
$ws.Range("K68").Value = 1833.3334
$ws.Range("M68").Value = -1084.3334

$ws.Range("H71").Value = 1833.3334
$ws.Range("I71").Value = 1833.3334
$ws.Range("K71").Value = 9166.666999999999
$ws.Range("M71").Value = -5422.666999999999

$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 11670.667
$ws.Range("I26").Value = 7506
$ws.Range("J26").Value = 20000
$ws.Range("K26").Value = 7506
$ws.Range("L26").Value = 20000
$ws.Range("M26").Value = -7213
$ws.Range("N26").Value = -20586

$ws.Range("H62").Value = 13280
$ws.Range("I62").Value = 9600
$ws.Range("K62").Value = 9600
$ws.Range("M62").Value = -8976

$ws.Range("H65").Value = 13280
$ws.Range("I65").Value = 9600
$ws.Range("K65").Value = 48000
$ws.Range("M65").Value = -44880

$ws.Range("H132").Value = 3790367.5
$ws.Range("I132").Value = 4506908.5
$ws.Range("J132").Value = 2936.2856
$ws.Range("K132").Value = 13520725.5
$ws.Range("L132").Value = 8808.856800000001
$ws.Range("M132").Value = -13518195.5
$ws.Range("N132").Value = -13868.8568
